# report of CHAMP measurement
#
# Adds the descriptive sentence about CHAMP optics to the (previously
# empty) second paragraph of the document, as six runs matching the
# author's original edit, each carrying an en-US language tag.

$d = $word.ActiveDocument

# The target is the empty paragraph that follows the "CHAMP Optics Test"
# title paragraph.
$target = $d.Paragraphs(2)
$insertionPoint = $target.Range.Start

$fragments = @(
    "CHAMP optics is the fore-optics of CHAI receiver array",
    ",",
    " which is used to ",
    "individually",
    " enlarge the ",
    "output beam waist of the Mixer block and adjust the "
)

foreach ($fragment in $fragments) {
    $run = $d.Range($insertionPoint, $insertionPoint)
    $run.InsertAfter($fragment)
    $run.LanguageID = "en-US"
    $insertionPoint = $run.End
}
